$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86-132 down to 87-133.
$ws.Rows.Item(86).Insert()

# Populate the new row 86 with the new week's data (weekly price update).
$ws.Cells.Item(86, 1).Value = 11
$ws.Cells.Item(86, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(86, 3).Value = "Bíobío"
$ws.Cells.Item(86, 4).Value = 44719
$ws.Cells.Item(86, 5).Value = 8
$ws.Cells.Item(86, 6).Value = 100112043
$ws.Cells.Item(86, 7).Value = "Pepino ensalada"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 100
$ws.Cells.Item(86, 11).Value = 20000
$ws.Cells.Item(86, 12).Value = 22000
$ws.Cells.Item(86, 13).Value = 21000
$ws.Cells.Item(86, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(86, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(86, 16).Value = 350
$ws.Cells.Item(86, 17).Value = 60
$ws.Cells.Item(86, 18).Value = "Hortaliza"

# Match the date-column style used elsewhere (numFmtId 165, style index 2 in the original file).
$ws.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
